# Updated symbol list on Sat Dec 31 05:16:20 UTC 2022 with GitHub Actions
#
# Applies the cell-level value updates to Sheet1 of the cryptos workbook:
#  - price (column D) refreshes for many rows
#  - three-way-ish re-shuffle of coin rows 10-19 (Coin/Link/Price/Volume columns)
#  - re-shuffle of rows 42/43 (KickToken <-> CEJI)
#  - a couple of "Volume(1h)" (column E) label tweaks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a number (so Excel would otherwise coerce
# the cell to a numeric type) get a leading apostrophe to force text, exactly
# like the source workbook stores these as text/inline strings.
$textUpdates = @(
    @{ Cell = "D2";  Value = "245.39" }
    @{ Cell = "D3";  Value = "25.41" }
    @{ Cell = "D5";  Value = "0.05574" }
    @{ Cell = "D6";  Value = "6.496" }
    @{ Cell = "D7";  Value = "3.019" }
    @{ Cell = "D8";  Value = "0.8182" }
    @{ Cell = "D9";  Value = "0.8463" }

    @{ Cell = "D10"; Value = "0.1340" }
    @{ Cell = "D11"; Value = "0.02876" }
    @{ Cell = "D12"; Value = "0.09377" }
    @{ Cell = "D13"; Value = "0.001512" }
    @{ Cell = "D14"; Value = "0.0005970" }
    @{ Cell = "D15"; Value = "0.006098" }
    @{ Cell = "D16"; Value = "3.497" }
    @{ Cell = "D17"; Value = "2.092" }
    @{ Cell = "D18"; Value = "0.3179" }
    @{ Cell = "D19"; Value = "0.06956" }

    @{ Cell = "D22"; Value = "3.750" }
    @{ Cell = "D23"; Value = "0.04724" }
    @{ Cell = "D25"; Value = "0.001250" }
    @{ Cell = "D26"; Value = "0.004636" }
    @{ Cell = "D27"; Value = "0.00009700" }
    @{ Cell = "D28"; Value = "0.0001390" }

    @{ Cell = "D40"; Value = "0.03663" }
    @{ Cell = "D41"; Value = "0.1354" }
    @{ Cell = "D42"; Value = "0.002660" }
    @{ Cell = "D43"; Value = "0.003379" }
    @{ Cell = "D44"; Value = "0.008303" }
    @{ Cell = "D45"; Value = "0.00005295" }
    @{ Cell = "D48"; Value = "0.002121" }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
}

# Plain text updates (coin names, links, volume labels) - no numeric coercion risk.
$plainUpdates = @(
    @{ Cell = "B10"; Value = "WazirX" }
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" }
    @{ Cell = "E10"; Value = "9WazirXWRX" }

    @{ Cell = "B11"; Value = "BitrueCoin" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "E11"; Value = "10BitrueCoinBTR" }

    @{ Cell = "B12"; Value = "BitMartToken" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Cell = "E12"; Value = "11BitMartTokenBMX" }

    @{ Cell = "B13"; Value = "BitForexToken" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Cell = "E13"; Value = "12BitForexTokenBF" }

    @{ Cell = "B14"; Value = "One" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" }
    @{ Cell = "E14"; Value = "13OneONE" }

    @{ Cell = "B15"; Value = "TigerCash" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "E15"; Value = "14TigerCashTCH" }

    @{ Cell = "B16"; Value = "LEO" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "E16"; Value = "15LEOLEO" }

    @{ Cell = "B17"; Value = "BTSEToken" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" }
    @{ Cell = "E17"; Value = "16BTSETokenBTSE" }

    @{ Cell = "B18"; Value = "BitpandaEcosystemToken" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best" }
    @{ Cell = "E18"; Value = "17BitpandaEcosystemTokenBEST" }

    @{ Cell = "B19"; Value = "MandalaExchangeToken" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Cell = "E19"; Value = "18MandalaExchangeTokenMDX" }

    @{ Cell = "E27"; Value = "26NitroExNTXBestin24h" }

    @{ Cell = "B42"; Value = "CEJI" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji" }
    @{ Cell = "E42"; Value = "41CEJICEJI" }

    @{ Cell = "B43"; Value = "KickToken" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick" }
    @{ Cell = "E43"; Value = "42KickTokenKICK" }
)

foreach ($u in $plainUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
